{"js": "// Commit: \"add words in June 7th\"\n//\n// Paragraph \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\" - the trailing \"\u65e5\u661f\u671f\" + \"\u4e8c\" runs get\n// coalesced into a single \"\u65e5\u661f\u671f\u4e8c\" run by the edit, and the following\n// weather paragraph (\"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7684\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\")\n// gains a new trailing sentence (\"\u3002\u4eca\u5929\u5929\u6c14\u4e0d\u9519\") in its own run, with\n// the paragraph-mark's direct formatting (the <w:pPr><w:rPr>) dropped.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"text\");\nawait context.sync();\n\n// --- Change 1: \"...\u65e5\u661f\u671f\" + \"\u4e8c\" -> \"...\u65e5\u661f\u671f\u4e8c\" (single run) ---------\n// Re-typing the full text over the existing (multi-run) match makes Word\n// coalesce the matched span into one run that uses the formatting of the\n// first run it replaced.\nconst dateMatches = body.search(\"\u65e5\u661f\u671f\u4e8c\", { matchCase: true });\nawait context.sync();\nif (dateMatches.items.length > 0) {\n  dateMatches.items[0].insertText(\"\u65e5\u661f\u671f\u4e8c\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2: append \"\u3002\u4eca\u5929\u5929\u6c14\u4e0d\u9519\" as its own run on the weather\n//     paragraph right after the June 7th date line, and drop that\n//     paragraph's stored paragraph-mark formatting. ------------------------\nbody.paragraphs.load(\"text\");\nawait context.sync();\n\nlet targetPara = null;\nfor (let i = 0; i < body.paragraphs.items.length; i++) {\n  const p = body.paragraphs.items[i];\n  if (p.text.indexOf(\"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7684\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\") !== -1) {\n    targetPara = p;\n    break;\n  }\n}\n\nif (targetPara) {\n  const whole = targetPara.getRange(Word.RangeLocation.whole);\n\n  // Read the paragraph's own OOXML so we keep its paraId/rsid* attributes\n  // and its existing run(s) exactly as-is, and only (a) drop the stored\n  // paragraph-mark formatting (<w:pPr>) and (b) append the new sentence as\n  // its own run (same rFonts hint as the run before it).\n  const existingOoxml = targetPara.getOoxml();\n  await context.sync();\n\n  const pMatch = existingOoxml.value.match(/<w:p\\b[^>]*>[\\s\\S]*?<\\/w:p>/);\n  let pXml = pMatch ? pMatch[0] : \"<w:p></w:p>\";\n\n  // Drop the paragraph-mark's own formatting block, if present.\n  pXml = pXml.replace(/<w:pPr>[\\s\\S]*?<\\/w:pPr>/, \"\");\n\n  // Append the new sentence as an additional run, just before </w:p>.\n  const newRun =\n    '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u3002\u4eca\u5929\u5929\u6c14\u4e0d\u9519</w:t></w:r>';\n  pXml = pXml.replace(/<\\/w:p>$/, newRun + \"</w:p>\");\n\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    '<w:body>' +\n    pXml +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  whole.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Commit: \"add words in June 7th\"\n#\n# Paragraph \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\" - the trailing \"\u65e5\u661f\u671f\" + \"\u4e8c\" runs get\n# coalesced into a single \"\u65e5\u661f\u671f\u4e8c\" run by the edit, and the following\n# weather paragraph (\"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7684\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\")\n# gains a new trailing sentence (\"\u3002\u4eca\u5929\u5929\u6c14\u4e0d\u9519\") in its own run, with\n# the paragraph-mark's direct formatting (the <w:pPr><w:rPr>) dropped.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"...\u65e5\u661f\u671f\" + \"\u4e8c\" -> \"...\u65e5\u661f\u671f\u4e8c\" (single run) -----------\n# A Find/Replace-All over the (multi-run) match makes Word coalesce the\n# matched span into a single run.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$found = $find.Execute(\n    \"\u65e5\u661f\u671f\u4e8c\",   # FindText\n    $false,       # MatchCase\n    $false,       # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    \"\u65e5\u661f\u671f\u4e8c\",   # ReplaceWith\n    2             # Replace (wdReplaceAll)\n)\n\n# --- Change 2: append \"\u3002\u4eca\u5929\u5929\u6c14\u4e0d\u9519\" as its own run on the weather\n#     paragraph right after the June 7th date line, and drop that\n#     paragraph's stored paragraph-mark formatting. --------------------------\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7684\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\") {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -ne $null) {\n    $r = $targetPara.Range\n\n    # Read the paragraph's own OOXML so we keep its paraId/rsid* attributes\n    # and its existing run(s) exactly as-is, and only (a) drop the stored\n    # paragraph-mark formatting (<w:pPr>) and (b) append the new sentence as\n    # its own run (same rFonts hint as the run before it).\n    $existingXml = $r.WordOpenXML\n\n    $pMatch = [regex]::Match($existingXml, '<w:p\\b[^>]*>[\\s\\S]*?</w:p>')\n    $pXml = $pMatch.Value\n\n    # Drop the paragraph-mark's own formatting block, if present.\n    $pXml = [regex]::Replace($pXml, '<w:pPr>[\\s\\S]*?</w:pPr>', '')\n\n    # Append the new sentence as an additional run, just before </w:p>.\n    $newRun = '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u3002\u4eca\u5929\u5929\u6c14\u4e0d\u9519</w:t></w:r>'\n    $pXml = $pXml -replace '</w:p>$', ($newRun + '</w:p>')\n\n    $ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n        '<w:body>' + $pXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $r.InsertXML($ooxml)\n}\n"}
